# Insert a new data row at row 505 (pushing the existing rows 505-611 down to
# 506-612) and populate the newly inserted row with its data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 505; this shifts rows 505..611
# down to 506..612 and keeps all their original values/formatting intact.
$ws.Rows.Item(505).Insert()

# Populate the newly blank row 505 with the new record's data.
$ws.Range("A505").Value = 5
$ws.Range("B505").Value = "Macroferia Regional de Talca"
$ws.Range("C505").Value = "Maule"
$ws.Range("D505").Value = 45211
$ws.Range("E505").Value = 7
$ws.Range("F505").Value = 100114013
$ws.Range("G505").Value = "Zanahoria"
$ws.Range("H505").Value = "Sin especificar"
$ws.Range("I505").Value = "Primera"
$ws.Range("J505").Value = 500
$ws.Range("K505").Value = 5000
$ws.Range("L505").Value = 5000
$ws.Range("M505").Value = 5000
$ws.Range("N505").Value = "$/saco 20 kilos"
$ws.Range("O505").Value = "Región de Ñuble"
$ws.Range("P505").Value = 250
$ws.Range("Q505").Value = 20
$ws.Range("R505").Value = "Hortaliza"
